$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 285
$ws.Range("I41").Value = 285
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 285
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 155
$ws.Range("N41").ClearContents()

$ws.Range("H53").Value = 863
$ws.Range("I53").Value = 858.55554
$ws.Range("J53").Value = 873
$ws.Range("K53").Value = 858.55554
$ws.Range("L53").Value = 873
$ws.Range("M53").Value = -221.55554
$ws.Range("N53").Value = -2147

$ws.Range("H76").Value = 3431.5454
$ws.Range("I76").Value = 3324.7
$ws.Range("K76").Value = 3324.7
$ws.Range("M76").Value = -3009.7

$ws.Range("H79").Value = 3431.5454
$ws.Range("I79").Value = 3324.7
$ws.Range("K79").Value = 3324.7
$ws.Range("M79").Value = -2232.7

$ws.Range("H86").Value = 4512.143
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 4964.1665
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 4964.1665
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -7210.1665

$ws.Range("H89").Value = 4512.143
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 4964.1665
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 24820.8325
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -36052.8325

$ws.Range("H92").Value = 1020.3
$ws.Range("I92").Value = 1020.3
$ws.Range("K92").Value = 1020.3
$ws.Range("M92").Value = 227.7

$ws.Range("H98").Value = 474.33334
$ws.Range("I98").Value = 474.33334
$ws.Range("K98").Value = 474.33334
$ws.Range("M98").Value = 1023.66666

$ws.Range("H113").Value = 19832.334
$ws.Range("I113").Value = 19832.334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 19832.334
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -16578.334
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 474.33334
$ws.Range("I122").Value = 474.33334
$ws.Range("K122").Value = 1423.00002
$ws.Range("M122").Value = 1026.99998

$ws.Range("H132").Value = 2050.7778
$ws.Range("I132").Value = 2050.7778
$ws.Range("K132").Value = 6152.3334
$ws.Range("M132").Value = -3622.3334

$ws.Range("H138").Value = 9000
$ws.Range("I138").Value = 9000
$ws.Range("K138").Value = 27000
$ws.Range("M138").Value = -21860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4279.2
$ws.Range("I32").Value = 3283.2942
$ws.Range("K32").Value = 3283.2942
$ws.Range("M32").Value = -2996.2942

$ws.Range("H61").Value = 1989
$ws.Range("I61").Value = 1989
$ws.Range("K61").Value = 1989
$ws.Range("M61").Value = -1777

$ws.Range("H97").Value = 1507.3334
$ws.Range("I97").Value = 1197.2632
$ws.Range("K97").Value = 1197.2632
$ws.Range("M97").Value = -701.2632000000001

$ws.Range("H110").Value = 919.5
$ws.Range("J110").Value = 850
$ws.Range("L110").Value = 850
$ws.Range("N110").Value = -4940

$ws.Range("H136").Value = 1989
$ws.Range("I136").Value = 1989
$ws.Range("K136").Value = 5967
$ws.Range("M136").Value = -3417

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 682.2174
$ws.Range("I107").Value = 681.4091
$ws.Range("K107").Value = 681.4091
$ws.Range("M107").Value = 1238.5909

$ws.Range("H134").Value = 2558.4285
$ws.Range("I134").Value = 2568.3333
$ws.Range("K134").Value = 7704.999899999999
$ws.Range("M134").Value = -5169.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1725.7727
$ws.Range("I7").Value = 897.8
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 897.8
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -784.8
$ws.Range("N7").Value = -3726

$ws.Range("H94").Value = 3240.25
$ws.Range("I94").Value = 3743.5
$ws.Range("J94").Value = 2737
$ws.Range("K94").Value = 3743.5
$ws.Range("L94").Value = 2737
$ws.Range("M94").Value = -3292.5
$ws.Range("N94").Value = -3639

$ws.Range("H105").Value = 16427.715
$ws.Range("I105").Value = 16427.715
$ws.Range("K105").Value = 16427.715
$ws.Range("M105").Value = -14680.715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1572.1333
$ws.Range("I102").Value = 1616
$ws.Range("J102").Value = 1451.5
$ws.Range("K102").Value = 1616
$ws.Range("L102").Value = 1451.5
$ws.Range("M102").Value = 6
$ws.Range("N102").Value = -4695.5

$ws.Range("H113").Value = 7540.375
$ws.Range("I113").Value = 5663.4
$ws.Range("K113").Value = 5663.4
$ws.Range("M113").Value = -3493.4

$ws.Range("H122").Value = 1899.8334
$ws.Range("I122").Value = 1975
$ws.Range("J122").Value = 1749.5
$ws.Range("K122").Value = 5925
$ws.Range("L122").Value = 5248.5
$ws.Range("M122").Value = -3475
$ws.Range("N122").Value = -10148.5

$ws.Range("H132").Value = 4442.4
$ws.Range("I132").Value = 4924.5
$ws.Range("K132").Value = 14773.5
$ws.Range("M132").Value = -12243.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 25000
$ws.Range("J36").Value = 25000
$ws.Range("L36").Value = 25000
$ws.Range("N36").Value = -26124

$ws.Range("H40").Value = 1780.5714
$ws.Range("I40").Value = 1660.6666
$ws.Range("K40").Value = 1660.6666
$ws.Range("M40").Value = -1524.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1599.5
$ws.Range("I62").Value = 1599.5
$ws.Range("K62").Value = 1599.5
$ws.Range("M62").Value = -975.5

$ws.Range("H65").Value = 1599.5
$ws.Range("I65").Value = 1599.5
$ws.Range("K65").Value = 7997.5
$ws.Range("M65").Value = -4877.5

$ws.Range("H122").Value = 733
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897

$ws.Range("H132").Value = 1747.5
$ws.Range("I132").Value = 1747.5
$ws.Range("K132").Value = 5242.5
$ws.Range("M132").Value = -2712.5

$ws.Range("H136").Value = 42628
$ws.Range("I136").Value = 46457.816
$ws.Range("J136").Value = 500
$ws.Range("K136").Value = 139373.448
$ws.Range("L136").Value = 1500
$ws.Range("M136").Value = -136823.448
$ws.Range("N136").Value = -6600
